$wb = $excel.ActiveWorkbook

# --- Rename "Tasks" sheet to "Remaining Tasks" ---
$ws1 = $wb.Worksheets.Item("02_03_24")
$ws3 = $wb.Worksheets.Item("Tasks")
$ws3.Name = "Remaining Tasks"

# --- Sheet "02_03_24": sprint updates ---
# (values assigned in this order so newly-appended shared strings land in
# the same order as the authored workbook)
$ws1.Range("E12").Value = "Placeholder "
$ws1.Range("D4").Value = "Started"
$ws1.Range("D12").Value = "Started"
$ws1.Range("D14").Value = "Started"
$ws1.Range("E7").Value = "Must rework shader and handle image processing class."
$ws1.Range("D13").Value = "Implemented"

# --- Sheet "Remaining Tasks": new Logger tasks + Application note ---
$ws3.Range("A62").Value = "Logger"
$ws3.Range("C62").Value = "Logger Class"
$ws3.Range("B62").Value = 1
$ws3.Range("C63").Value = "Updated all error messages to logger"
$ws3.Range("B63").Value = 2
$ws3.Range("E28").Value = "Needs to be singleton."

$ws1.Range("E2").Value = "Need to be a singleton."

# --- Restore selections to match the final saved view state ---
$ws3.Range("E28").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("E14").Select() | Out-Null
